$wb = $excel.ActiveWorkbook

# This workbook is refreshed on a schedule from live Universalis market-board
# data for the Excalibur data-center. Each worksheet tab (one per crafting job)
# lists FFXIV leves together with the current average NQ/HQ marketboard prices
# and the resulting leve-turn-in profit. Columns:
#   H currentAveragePrice    I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ
#   M LeveProfitNQ           N LeveProfitHQ
# A blank M/N means that quality tier had no sellable market data that refresh,
# so the profit cell is left empty (mirrors the upstream exporter, which omits
# the cell entirely rather than writing a zero/placeholder).

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 782.3333
$ws.Range("I4").Value = 165.57143
$ws.Range("J4").Value = 1645.8
$ws.Range("K4").Value = 165.57143
$ws.Range("L4").Value = 1645.8
$ws.Range("M4").Value = -51.57142999999999
$ws.Range("N4").Value = -1873.8
$ws.Range("H39").Value = 1799.3334
$ws.Range("I39").Value = 199.5
$ws.Range("J39").Value = 4999
$ws.Range("K39").Value = 598.5
$ws.Range("L39").Value = 14997
$ws.Range("M39").Value = -302.5
$ws.Range("N39").Value = -15589
$ws.Range("H132").Value = 54016.754
$ws.Range("I132").Value = 55598.953
$ws.Range("K132").Value = 166796.859
$ws.Range("M132").Value = -164266.859
$ws.Range("H137").Value = 4162.3335
$ws.Range("I137").Value = 3300.5908
$ws.Range("J137").Value = 5885.8184
$ws.Range("K137").Value = 9901.7724
$ws.Range("L137").Value = 17657.4552
$ws.Range("M137").Value = -7351.7724
$ws.Range("N137").Value = -22757.4552
$ws.Range("H138").Value = 4600.968
$ws.Range("J138").Value = 4401.0347
$ws.Range("L138").Value = 13203.1041
$ws.Range("N138").Value = -23483.1041

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6099262.5
$ws.Range("I32").Value = 7353674.5
$ws.Range("J32").Value = 6406
$ws.Range("K32").Value = 7353674.5
$ws.Range("L32").Value = 6406
$ws.Range("M32").Value = -7353387.5
$ws.Range("N32").Value = -6980
$ws.Range("H61").Value = 2900.2222
$ws.Range("I61").Value = 2705.5938
$ws.Range("J61").Value = 3379.3076
$ws.Range("K61").Value = 2705.5938
$ws.Range("L61").Value = 3379.3076
$ws.Range("M61").Value = -2493.5938
$ws.Range("N61").Value = -3803.3076
$ws.Range("H97").Value = 1127.1515
$ws.Range("I97").Value = 1107.5927
$ws.Range("K97").Value = 1107.5927
$ws.Range("M97").Value = -611.5926999999999
$ws.Range("H108").Value = 80684
$ws.Range("J108").Value = 80684
$ws.Range("L108").Value = 80684
$ws.Range("N108").Value = -88364
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = ""
$ws.Range("H132").Value = 2135.1384
$ws.Range("I132").Value = 1959
$ws.Range("J132").Value = 2999.818
$ws.Range("K132").Value = 5877
$ws.Range("L132").Value = 8999.454000000002
$ws.Range("M132").Value = -3347
$ws.Range("N132").Value = -14059.454
$ws.Range("H136").Value = 2900.2222
$ws.Range("I136").Value = 2705.5938
$ws.Range("J136").Value = 3379.3076
$ws.Range("K136").Value = 8116.7814
$ws.Range("L136").Value = 10137.9228
$ws.Range("M136").Value = -5566.7814
$ws.Range("N136").Value = -15237.9228

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2494.3333
$ws.Range("I86").Value = 1992.7858
$ws.Range("J86").Value = 4249.75
$ws.Range("K86").Value = 1992.7858
$ws.Range("L86").Value = 4249.75
$ws.Range("M86").Value = -869.7858000000001
$ws.Range("N86").Value = -6495.75
$ws.Range("H89").Value = 2494.3333
$ws.Range("I89").Value = 1992.7858
$ws.Range("J89").Value = 4249.75
$ws.Range("K89").Value = 9963.929
$ws.Range("L89").Value = 21248.75
$ws.Range("M89").Value = -4347.929
$ws.Range("N89").Value = -32480.75
$ws.Range("H99").Value = 11699.77
$ws.Range("I99").Value = 5734.6113
$ws.Range("J99").Value = 25121.375
$ws.Range("K99").Value = 5734.6113
$ws.Range("L99").Value = 25121.375
$ws.Range("M99").Value = -4236.6113
$ws.Range("N99").Value = -28117.375
$ws.Range("H105").Value = 1338.931
$ws.Range("I105").Value = 1321.64
$ws.Range("J105").Value = 1447
$ws.Range("K105").Value = 1321.64
$ws.Range("L105").Value = 1447
$ws.Range("M105").Value = 425.3599999999999
$ws.Range("N105").Value = -4941
$ws.Range("H107").Value = 1592.3914
$ws.Range("I107").Value = 1708.9445
$ws.Range("K107").Value = 1708.9445
$ws.Range("M107").Value = 211.0554999999999
$ws.Range("H108").Value = 100684
$ws.Range("J108").Value = 100684
$ws.Range("L108").Value = 100684
$ws.Range("N108").Value = -108364
$ws.Range("H134").Value = 2228.6743
$ws.Range("I134").Value = 1157.8108
$ws.Range("J134").Value = 8832.333000000001
$ws.Range("K134").Value = 3473.4324
$ws.Range("L134").Value = 26496.999
$ws.Range("M134").Value = -938.4323999999997
$ws.Range("N134").Value = -31566.999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 74701.22
$ws.Range("I31").Value = 96048.94
$ws.Range("K31").Value = 96048.94
$ws.Range("M31").Value = -95753.94
$ws.Range("H34").Value = 74701.22
$ws.Range("I34").Value = 96048.94
$ws.Range("K34").Value = 96048.94
$ws.Range("M34").Value = -95846.94
$ws.Range("H58").Value = 1895.3024
$ws.Range("I58").Value = 1243.9667
$ws.Range("J58").Value = 3398.3845
$ws.Range("K58").Value = 1243.9667
$ws.Range("L58").Value = 3398.3845
$ws.Range("M58").Value = -1040.9667
$ws.Range("N58").Value = -3804.3845
$ws.Range("H62").Value = 3199.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 3199.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H107").Value = 755.0625
$ws.Range("I107").Value = 650.7
$ws.Range("J107").Value = 929
$ws.Range("K107").Value = 650.7
$ws.Range("L107").Value = 929
$ws.Range("M107").Value = 1269.3
$ws.Range("N107").Value = -4769
$ws.Range("H132").Value = 8622301
$ws.Range("I132").Value = 1495.909
$ws.Range("J132").Value = 35716260
$ws.Range("K132").Value = 4487.727000000001
$ws.Range("L132").Value = 107148780
$ws.Range("M132").Value = -1957.727000000001
$ws.Range("N132").Value = -107153840
$ws.Range("H136").Value = 1895.3024
$ws.Range("I136").Value = 1243.9667
$ws.Range("J136").Value = 3398.3845
$ws.Range("K136").Value = 3731.9001
$ws.Range("L136").Value = 10195.1535
$ws.Range("M136").Value = -1181.9001
$ws.Range("N136").Value = -15295.1535

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 980.2857
$ws.Range("I5").Value = 978.8333
$ws.Range("J5").Value = 989
$ws.Range("K5").Value = 2936.4999
$ws.Range("L5").Value = 2967
$ws.Range("M5").Value = -2824.4999
$ws.Range("N5").Value = -3191
$ws.Range("H20").Value = 4553.1665
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("H33").Value = 176.66667
$ws.Range("I33").Value = 176.66667
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1060.00002
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -777.0000199999999
$ws.Range("N33").Value = ""
$ws.Range("H39").Value = 7999
$ws.Range("I39").Value = 7999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 23997
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -23703
$ws.Range("N39").Value = ""
$ws.Range("H132").Value = 3825
$ws.Range("I132").Value = 3825
$ws.Range("K132").Value = 34425
$ws.Range("M132").Value = -31895
$ws.Range("H135").Value = 980.2857
$ws.Range("I135").Value = 978.8333
$ws.Range("J135").Value = 989
$ws.Range("K135").Value = 8809.4997
$ws.Range("L135").Value = 8901
$ws.Range("M135").Value = -6274.4997
$ws.Range("N135").Value = -13971

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 29399.2
$ws.Range("J55").Value = 29399.2
$ws.Range("L55").Value = 29399.2
$ws.Range("N55").Value = -30053.2
$ws.Range("H62").Value = 76723.336
$ws.Range("H65").Value = 76723.336
$ws.Range("H107").Value = 2062.3635
$ws.Range("J107").Value = 1715.5714
$ws.Range("L107").Value = 1715.5714
$ws.Range("N107").Value = -5555.5714
$ws.Range("H108").Value = 96102.60000000001
$ws.Range("J108").Value = 96102.60000000001
$ws.Range("L108").Value = 96102.60000000001
$ws.Range("N108").Value = -103782.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 728.4194
$ws.Range("I16").Value = 686.95654
$ws.Range("J16").Value = 847.625
$ws.Range("K16").Value = 686.95654
$ws.Range("L16").Value = 847.625
$ws.Range("M16").Value = -516.95654
$ws.Range("N16").Value = -1187.625
$ws.Range("H46").Value = 2502.4375
$ws.Range("J46").Value = 2225.3076
$ws.Range("L46").Value = 2225.3076
$ws.Range("N46").Value = -2601.3076
$ws.Range("H55").Value = 751.4737
$ws.Range("I55").Value = 200.85715
$ws.Range("J55").Value = 1072.6666
$ws.Range("K55").Value = 200.85715
$ws.Range("L55").Value = 1072.6666
$ws.Range("M55").Value = -27.85714999999999
$ws.Range("N55").Value = -1418.6666
$ws.Range("H100").Value = 5375
$ws.Range("I100").Value = 1819.0769
$ws.Range("J100").Value = 51602
$ws.Range("K100").Value = 1819.0769
$ws.Range("L100").Value = 51602
$ws.Range("M100").Value = -1278.0769
$ws.Range("N100").Value = -52684
$ws.Range("H122").Value = 129664.625
$ws.Range("J122").Value = 171302.33
$ws.Range("L122").Value = 513906.99
$ws.Range("N122").Value = -518806.99
$ws.Range("H132").Value = 2082.8333
$ws.Range("I132").Value = 1835.9615
$ws.Range("K132").Value = 5507.8845
$ws.Range("M132").Value = -2977.8845

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 54999.668
$ws.Range("I54").Value = 37500
$ws.Range("J54").Value = 89999
$ws.Range("K54").Value = 37500
$ws.Range("L54").Value = 89999
$ws.Range("M54").Value = -36980
$ws.Range("N54").Value = -91039
$ws.Range("H107").Value = 3185.6667
$ws.Range("I107").Value = 2243.6875
$ws.Range("J107").Value = 5069.625
$ws.Range("K107").Value = 6731.0625
$ws.Range("L107").Value = 15208.875
$ws.Range("M107").Value = -4811.0625
$ws.Range("N107").Value = -19048.875
